$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.273.39'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '1.551.78'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.09'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0582'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '1.774.32'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").Value = '1.551.65'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.62'
$ws.Range("E15").Value = '  -0.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.507'
$ws.Range("E16").Value = '  -2.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.26'
$ws.Range("E17").Value = '  -2.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.28'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.29'
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").Value = '0.0₃0670'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("E23").Value = '  -2.95%  '
$ws.Range("E24").Value = '  -4.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.37'
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.73'
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.22'
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0467'
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.03'
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").Value = '1.382.61'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("E36").Value = '  -3.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.58'
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.510'
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.91'
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.910'
$ws.Range("E47").Value = '  -6.18%  '
$ws.Range("D48").Value = '1.687.00'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.07'
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0103'
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '41.51'
$ws.Range("E51").Value = '  +8.05%  '
